# DateBox and AdvAPI final demo modifications
# Appends a new daily work-report block (rows 173-179, for 2025-02-10)
# to Sheet1, mirroring the content/formatting of the preceding blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Populate the new rows' values/formulas first --------------------
# (Formatting copies are applied afterwards so that the live formula
# engine recalculates against the final values.)
#
# NOTE: the two brand-new shared strings ("Reconsile Revision" and
# "Overview: Data Layer & Data Grid") must be written in this exact
# order so they land at shared-string indices 50 and 51 respectively,
# matching the workbook being reproduced - hence row 177's text is
# assigned before row 174's below.

# Row 173: date + "Domm"
$ws.Range("A173").Value = 45698
$ws.Range("B173").Value = "Domm"
$ws.Range("D173").Value = 0.25

# Row 177: Study / Reconsile Revision (text set early - see note above)
$ws.Range("B177").Value = "Study"
$ws.Range("C177").Value = "Reconsile Revision"
$ws.Range("D177").Value = 2

# Row 174: Meeting / Overview: Data Layer & Data Grid
$ws.Range("B174").Value = "Meeting"
$ws.Range("C174").Value = "Overview: Data Layer & Data Grid"
$ws.Range("D174").Value = 0.25

# Row 175: Reconsile (no entry in column B, matching source data)
$ws.Range("C175").Value = "Reconsile"
$ws.Range("D175").Value = 1

# Row 176: General Discussion
$ws.Range("C176").Value = "General Discussion"
$ws.Range("D176").Value = 0.25

# Row 178: Editors - Overview
$ws.Range("C178").Value = "Editors – Overview"
$ws.Range("D178").Value = 4.25

# Row 179: Total
$ws.Range("B179").Value = "Total"
$ws.Range("D179").Formula = "=SUM(D172:D178)"

# --- 2) Copy cell formatting down from the previous block ---------------
# Row 167 carries the "date style" formatting for column A (style 2),
# shared by rows 173-175 here; row 170 carries the "detail/total row"
# formatting (style 3 for column A), shared by rows 176-179 here.
$ws.Range("A167:D167").Copy()
$ws.Range("A173:D175").PasteSpecial(-4122)

$ws.Range("A170:D170").Copy()
$ws.Range("A176:D179").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 175 has no entry in column B at all (not even a blank styled
# cell) in the source data, so drop the formatting the tiled paste put
# there.
$ws.Range("B175").Clear()
$ws.Range("C175").Value = "Reconsile"
$ws.Range("D175").Value = 1

# --- 3) Update the view state (scroll position / selection) ------------
$ws.Application.Goto($ws.Range("A147"), $true)
[void]$ws.Range("U158").Select()

Write-Output "Applied DailyWorkReport update (rows 173-179)."
